$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7842.75
$ws.Range("I18").Value = 9333
$ws.Range("J18").Value = 3372
$ws.Range("K18").Value = 9333
$ws.Range("L18").Value = 3372
$ws.Range("M18").Value = -9049
$ws.Range("N18").Value = -3940

$ws.Range("H28").Value = 7429.4
$ws.Range("I28").Value = 25399.5
$ws.Range("J28").Value = 2936.875
$ws.Range("K28").Value = 25399.5
$ws.Range("L28").Value = 2936.875
$ws.Range("M28").Value = -24914.5
$ws.Range("N28").Value = -3906.875

$ws.Range("H69").Value = 12179
$ws.Range("J69").Value = 16947.5
$ws.Range("L69").Value = 50842.5
$ws.Range("N69").Value = -52590.5

$ws.Range("H72").Value = 12179
$ws.Range("J72").Value = 16947.5
$ws.Range("L72").Value = 152527.5
$ws.Range("N72").Value = -161263.5

$ws.Range("H86").Value = 59062.375
$ws.Range("I86").Value = 3125
$ws.Range("K86").Value = 3125
$ws.Range("M86").Value = -2002

$ws.Range("H89").Value = 59062.375
$ws.Range("I89").Value = 3125
$ws.Range("K89").Value = 15625
$ws.Range("M89").Value = -10009

$ws.Range("H98").Value = 37554.723
$ws.Range("I98").Value = 40499.285
$ws.Range("K98").Value = 40499.285
$ws.Range("M98").Value = -39001.285

$ws.Range("H122").Value = 37554.723
$ws.Range("I122").Value = 40499.285
$ws.Range("K122").Value = 121497.855
$ws.Range("M122").Value = -119047.855

$ws.Range("H136").Value = 65499.5
$ws.Range("J136").Value = 65499.5
$ws.Range("L136").Value = 65499.5
$ws.Range("N136").Value = -75699.5

$ws.Range("H137").Value = 8413.9375
$ws.Range("I137").Value = 9137.357
$ws.Range("K137").Value = 27412.071
$ws.Range("M137").Value = -24862.071

$ws.Range("H138").Value = 3738.1162
$ws.Range("I138").Value = 1144.4445
$ws.Range("J138").Value = 4424.6763
$ws.Range("K138").Value = 3433.3335
$ws.Range("L138").Value = 13274.0289
$ws.Range("M138").Value = 1706.6665
$ws.Range("N138").Value = -23554.0289

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1973.3
$ws.Range("I32").Value = 1972.1875
$ws.Range("K32").Value = 1972.1875
$ws.Range("M32").Value = -1685.1875

$ws.Range("H95").Value = 29999
$ws.Range("J95").Value = 29999
$ws.Range("L95").Value = 29999
$ws.Range("N95").Value = -35491

$ws.Range("H97").Value = 16105.154
$ws.Range("I97").Value = 7858.2856
$ws.Range("K97").Value = 7858.2856
$ws.Range("M97").Value = -7362.2856

$ws.Range("H122").Value = 471935.72
$ws.Range("I122").Value = 3430.3572
$ws.Range("K122").Value = 10291.0716
$ws.Range("M122").Value = -7841.071599999999

$ws.Range("H132").Value = 4086.2092
$ws.Range("I132").Value = 3332.7576
$ws.Range("J132").Value = 6572.6
$ws.Range("K132").Value = 9998.272799999999
$ws.Range("L132").Value = 19717.8
$ws.Range("M132").Value = -7468.272799999999
$ws.Range("N132").Value = -24777.8

$ws.Range("H141").Value = 85976.336
$ws.Range("J141").Value = 85976.336
$ws.Range("L141").Value = 85976.336
$ws.Range("N141").Value = -96336.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 15478.615
$ws.Range("I99").Value = 15478.615
$ws.Range("K99").Value = 15478.615
$ws.Range("M99").Value = -13980.615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3444.9092
$ws.Range("I31").Value = 1165.8334
$ws.Range("K31").Value = 1165.8334
$ws.Range("M31").Value = -870.8334

$ws.Range("H34").Value = 3444.9092
$ws.Range("I34").Value = 1165.8334
$ws.Range("K34").Value = 1165.8334
$ws.Range("M34").Value = -963.8334

$ws.Range("H58").Value = 1628.25
$ws.Range("J58").Value = 8000
$ws.Range("L58").Value = 8000
$ws.Range("N58").Value = -8406

$ws.Range("H97").Value = 79399.75
$ws.Range("J97").Value = 79399.75
$ws.Range("L97").Value = 79399.75
$ws.Range("N97").Value = -81381.75

$ws.Range("H109").Value = 76047.336
$ws.Range("J109").Value = 76047.336
$ws.Range("L109").Value = 76047.336
$ws.Range("N109").Value = -78127.336

$ws.Range("H122").Value = 2193.3333
$ws.Range("I122").Value = 2200
$ws.Range("J122").Value = 2190
$ws.Range("K122").Value = 6600
$ws.Range("L122").Value = 6570
$ws.Range("M122").Value = -4150
$ws.Range("N122").Value = -11470

$ws.Range("H136").Value = 1628.25
$ws.Range("J136").Value = 8000
$ws.Range("L136").Value = 24000
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 44139788
$ws.Range("I4").Value = 34839868
$ws.Range("K4").Value = 104519604
$ws.Range("M4").Value = -104519492

$ws.Range("H5").Value = 456810.12
$ws.Range("I5").Value = 2147.5386
$ws.Range("J5").Value = 1113545
$ws.Range("K5").Value = 6442.6158
$ws.Range("L5").Value = 3340635
$ws.Range("M5").Value = -6330.6158
$ws.Range("N5").Value = -3340859

$ws.Range("H118").Value = 8750

$ws.Range("H135").Value = 456810.12
$ws.Range("I135").Value = 2147.5386
$ws.Range("J135").Value = 1113545
$ws.Range("K135").Value = 19327.8474
$ws.Range("L135").Value = 10021905
$ws.Range("M135").Value = -16792.8474
$ws.Range("N135").Value = -10026975

$ws.Range("H140").Value = 13028.6
$ws.Range("I140").Value = 14417.692
$ws.Range("J140").Value = 3999.5
$ws.Range("K140").Value = 43253.076
$ws.Range("L140").Value = 11998.5
$ws.Range("M140").Value = -38073.076
$ws.Range("N140").Value = -22358.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H70").Value = 6108.5186
$ws.Range("I70").Value = 5752.4165
$ws.Range("J70").Value = 6393.4
$ws.Range("K70").Value = 5752.4165
$ws.Range("L70").Value = 6393.4
$ws.Range("M70").Value = -5482.4165
$ws.Range("N70").Value = -6933.4

$ws.Range("H73").Value = 6108.5186
$ws.Range("I73").Value = 5752.4165
$ws.Range("J73").Value = 6393.4
$ws.Range("K73").Value = 5752.4165
$ws.Range("L73").Value = 6393.4
$ws.Range("M73").Value = -4816.4165
$ws.Range("N73").Value = -8265.4

$ws.Range("H113").Value = 2000.3334
$ws.Range("I113").Value = 1600.4
$ws.Range("K113").Value = 1600.4
$ws.Range("M113").Value = 569.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3064.1765
$ws.Range("J22").Value = 2207.125
$ws.Range("L22").Value = 2207.125
$ws.Range("N22").Value = -2797.125

$ws.Range("H27").Value = 3064.1765
$ws.Range("J27").Value = 2207.125
$ws.Range("L27").Value = 2207.125
$ws.Range("N27").Value = -2421.125

$ws.Range("H61").Value = 4725.516
$ws.Range("I61").Value = 3145.923
$ws.Range("K61").Value = 3145.923
$ws.Range("M61").Value = -2943.923

$ws.Range("H109").Value = 60000
$ws.Range("I109").Value = 60000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 60000
$ws.Range("L109").ClearContents()
$ws.Range("M109").Value = -58613
$ws.Range("N109").Value = 0

$ws.Range("H113").Value = 4725.516
$ws.Range("I113").Value = 3145.923
$ws.Range("K113").Value = 3145.923
$ws.Range("M113").Value = -975.9229999999998

$ws.Range("H131").Value = 84500
$ws.Range("J131").Value = 84500
$ws.Range("L131").Value = 84500
$ws.Range("N131").Value = -94580

$ws.Range("H132").Value = 455128.47
$ws.Range("I132").Value = 786370.75
$ws.Range("K132").Value = 2359112.25
$ws.Range("M132").Value = -2356582.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 12506300
$ws.Range("J5").Value = 12600
$ws.Range("L5").Value = 12600
$ws.Range("N5").Value = -12824

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0

$ws.Range("H113").Value = 2894.7693
$ws.Range("I113").Value = 1469.0834
$ws.Range("J113").Value = 20003
$ws.Range("K113").Value = 4407.2502
$ws.Range("L113").Value = 60009
$ws.Range("M113").Value = -2237.2502
$ws.Range("N113").Value = -64349
